$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Part A: consolidate runs that now carry identical (unchanged) visible text
# into a single run per paragraph. We do this by re-finding the already
# correct text and "replacing" it with itself inside the Range of each
# specific paragraph (found by index) - this makes Word coalesce the runs
# without altering any visible content.
# ---------------------------------------------------------------------------

function Merge-ParagraphText($paraIndex, $exactText) {
    $p = $d.Paragraphs.Item($paraIndex)
    $r = $p.Range
    $r.Find.Execute($exactText, $true, $false, $false, $false, $false, $true, 1, $false, $exactText, 2) | Out-Null
}

Merge-ParagraphText 95 "5. ik denk Konijn! Omdat er geen break op maandag staat "
Merge-ParagraphText 98 "8. Konijn!"
Merge-ParagraphText 99 "9. Konijn!"
Merge-ParagraphText 100 "10. Konijn!"
Merge-ParagraphText 101 "11. vos!"
Merge-ParagraphText 102 "12. Konijn!"

# ---------------------------------------------------------------------------
# Part B: append the new "Part 2" / "Part 3" answers at the end of the body.
# ---------------------------------------------------------------------------

function Add-Paragraph($text, [bool]$englishUK) {
    $last = $d.Paragraphs.Last
    $last.Range.InsertParagraphAfter() | Out-Null
    $newPara = $d.Paragraphs.Last
    $nr = $newPara.Range
    $nr.Text = $text
    if ($englishUK) {
        $d.Paragraphs.Last.Range.LanguageID = "en-GB"
    }
}

Add-Paragraph "Part 2 – 4" $false
Add-Paragraph "1, great omdat het persoon 18 jaar is en geboren is in 1988" $false
Add-Paragraph "Part 3 – 1" $true
Add-Paragraph "3, String" $true
Add-Paragraph "4. interger" $true
Add-Paragraph "5. String" $true
Add-Paragraph "6. Boolean" $true
Add-Paragraph "7. String" $true
Add-Paragraph " " $true

Write-Host "Edit complete"
